$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1850746268656716
$ws.Range("C2").Value = 0.5611940298507463
$ws.Range("J2").Value = 0.0208955223880597
$ws.Range("P2").Value = 0.1373134328358209
$ws.Range("S2").Value = 0.09552238805970149
$ws.Range("B3").Value = 0.0154639175257732
$ws.Range("C3").Value = 0.03092783505154639
$ws.Range("J3").Value = 0.04123711340206185
$ws.Range("P3").Value = 0.7731958762886598
$ws.Range("S3").Value = 0.1391752577319588
$ws.Range("J4").Value = 0.1458333333333333
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.1875
$ws.Range("B6").Value = 0.09236947791164658
$ws.Range("D6").Value = 0.004016064257028112
$ws.Range("F6").Value = 0.07630522088353414
$ws.Range("J6").Value = 0.2008032128514056
$ws.Range("O6").Value = 0.01204819277108434
$ws.Range("Q6").Value = 0.1646586345381526
$ws.Range("R6").Value = 0.06024096385542169
$ws.Range("S6").Value = 0.3895582329317269
$ws.Range("B7").Value = 0.1213872832369942
$ws.Range("D7").Value = 0.01734104046242774
$ws.Range("F7").Value = 0.05780346820809248
$ws.Range("J7").Value = 0.161849710982659
$ws.Range("O7").Value = 0.04624277456647399
$ws.Range("Q7").Value = 0.1849710982658959
$ws.Range("R7").Value = 0.06936416184971098
$ws.Range("S7").Value = 0.3410404624277457
$ws.Range("B8").Value = 0.1044776119402985
$ws.Range("D8").Value = 0.02985074626865672
$ws.Range("F8").Value = 0.07213930348258707
$ws.Range("J8").Value = 0.1268656716417911
$ws.Range("O8").Value = 0.02238805970149254
$ws.Range("Q8").Value = 0.1691542288557214
$ws.Range("R8").Value = 0.1044776119402985
$ws.Range("S8").Value = 0.3706467661691542
$ws.Range("B9").Value = 0.08121827411167512
$ws.Range("D9").Value = 0.02030456852791878
$ws.Range("F9").Value = 0.1065989847715736
$ws.Range("J9").Value = 0.1065989847715736
$ws.Range("O9").Value = 0.04060913705583756
$ws.Range("Q9").Value = 0.1573604060913706
$ws.Range("R9").Value = 0.09137055837563451
$ws.Range("S9").Value = 0.3959390862944163
$ws.Range("B10").Value = 0.1230425055928412
$ws.Range("D10").Value = 0.02162565249813572
$ws.Range("E10").Value = 0.0007457121551081282
$ws.Range("F10").Value = 0.07979120059656973
$ws.Range("J10").Value = 0.116331096196868
$ws.Range("O10").Value = 0.01864280387770321
$ws.Range("Q10").Value = 0.2155108128262491
$ws.Range("R10").Value = 0.07606263982102908
$ws.Range("S10").Value = 0.3482475764354959
$ws.Range("G11").Value = 0.1346801346801347
$ws.Range("J11").Value = 0.1144781144781145
$ws.Range("K11").Value = 0.2154882154882155
$ws.Range("L11").Value = 0.5252525252525253
$ws.Range("S11").Value = 0.0101010101010101
$ws.Range("G12").Value = 0.7098765432098766
$ws.Range("J12").Value = 0.228395061728395
$ws.Range("K12").Value = 0.006172839506172839
$ws.Range("L12").Value = 0.0308641975308642
$ws.Range("S12").Value = 0.02469135802469136
$ws.Range("G13").Value = 0.7428571428571429
$ws.Range("J13").Value = 0.2571428571428571
$ws.Range("F15").Value = 0.01838235294117647
$ws.Range("H15").Value = 0.125
$ws.Range("I15").Value = 0.05147058823529412
$ws.Range("J15").Value = 0.3860294117647059
$ws.Range("K15").Value = 0.06985294117647059
$ws.Range("M15").Value = 0.01102941176470588
$ws.Range("O15").Value = 0.07352941176470588
$ws.Range("S15").Value = 0.2647058823529412
$ws.Range("F16").Value = 0.01345291479820628
$ws.Range("H16").Value = 0.1210762331838565
$ws.Range("I16").Value = 0.05829596412556054
$ws.Range("J16").Value = 0.4618834080717489
$ws.Range("K16").Value = 0.09417040358744394
$ws.Range("M16").Value = 0.02242152466367713
$ws.Range("O16").Value = 0.07623318385650224
$ws.Range("S16").Value = 0.1524663677130045
$ws.Range("F17").Value = 0.0196078431372549
$ws.Range("H17").Value = 0.1699346405228758
$ws.Range("I17").Value = 0.08932461873638345
$ws.Range("J17").Value = 0.4030501089324618
$ws.Range("K17").Value = 0.1154684095860566
$ws.Range("M17").Value = 0.01742919389978214
$ws.Range("N17").Value = 0.002178649237472767
$ws.Range("O17").Value = 0.08061002178649238
$ws.Range("S17").Value = 0.10239651416122
$ws.Range("F18").Value = 0.01612903225806452
$ws.Range("H18").Value = 0.1881720430107527
$ws.Range("I18").Value = 0.05913978494623656
$ws.Range("J18").Value = 0.4623655913978494
$ws.Range("K18").Value = 0.07526881720430108
$ws.Range("M18").Value = 0.01075268817204301
$ws.Range("N18").Value = 0.005376344086021506
$ws.Range("O18").Value = 0.08602150537634409
$ws.Range("S18").Value = 0.09677419354838709
$ws.Range("F19").Value = 0.01312551271534044
$ws.Range("H19").Value = 0.1862182116488925
$ws.Range("I19").Value = 0.09515996718621821
$ws.Range("J19").Value = 0.3904840032813782
$ws.Range("K19").Value = 0.0992616899097621
$ws.Range("M19").Value = 0.014766201804758
$ws.Range("O19").Value = 0.08285479901558655
$ws.Range("S19").Value = 0.118129614438064
